$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Re-theme the title + header fonts: both collapse onto a single bold,
#    white (FFFFFF), default-size (11) font instead of the previous pair of
#    separate "bold 14pt" (title) and "bold, default color" (header) fonts.
# ---------------------------------------------------------------------------
$white = 16777215   # RGB(255,255,255)

foreach ($ws in $wb.Worksheets) {
    # Title cell (row 1, col A)
    $title = $ws.Range("A1")
    $title.Font.Bold = $true
    $title.Font.Size = 11
    $title.Font.Color = $white

    # Header row (row 2) spans the sheet's used width
    $usedRange = $ws.UsedRange
    $lastCol = $usedRange.Columns.Count
    $headerRange = $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item(2, $lastCol))
    $headerRange.Font.Bold = $true
    $headerRange.Font.Size = 11
    $headerRange.Font.Color = $white
}

# ---------------------------------------------------------------------------
# 2) "Training Dashboard" sheet: refresh PERIOD TO EXPIRE (H) and
#    LAST UPDATE (I) after recalculating against a newer "last update" date.
# ---------------------------------------------------------------------------
$wsTraining = $wb.Worksheets.Item("Training Dashboard")

$periodToExpire = @{
    3 = 677;  4 = 565;  5 = 566;  6 = 564;  7 = 566;  8 = 566;
    9 = 565;  10 = 564; 11 = 697; 12 = 566; 13 = 478; 14 = 378;
    15 = 435; 16 = 482; 17 = 423; 18 = 200; 19 = 200; 20 = 312;
    21 = 312; 22 = 312; 23 = 333; 24 = 333
}

for ($row = 3; $row -le 24; $row++) {
    $wsTraining.Cells.Item($row, 8).Value = $periodToExpire[$row]
    # Leading apostrophe forces literal text so "16-Sep-2025" is stored as a
    # string (matching the original inline-string cells) instead of being
    # auto-converted into a serial date value.
    $wsTraining.Cells.Item($row, 9).Value = "'16-Sep-2025"
}

# ---------------------------------------------------------------------------
# 3) "Exam Dashboard" sheet: widen the COMMENTS column and refresh its
#    wording now that dates are being validated explicitly.
# ---------------------------------------------------------------------------
$wsExam = $wb.Worksheets.Item("Exam Dashboard")

$wsExam.Columns.Item(5).ColumnWidth = 14.17

for ($row = 3; $row -le 11; $row++) {
    $wsExam.Cells.Item($row, 5).Value = "date is valid"
}
